$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(1, 8).Value = "AnswerIsCorrect"
$ws.Range("A2:G10").EntireRow.Delete()
$ws.Range("A11:H16").RowHeight = 15.75
$ws.Range("A17:A20").RowHeight = 15.75
